$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: find the first paragraph whose text starts with $needle.
# ------------------------------------------------------------------
function Find-ParagraphByText($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($needle)) {
            return $p
        }
    }
    return $null
}

# ------------------------------------------------------------------
# 1. Insert the three replacement bullet paragraphs right after the
#    "Churn and Tenure..." bullet (a numId=8 list item), WHILE the
#    "Advanced Analysis" section (which still follows it) is in
#    place - InsertXML on a placeholder paragraph that is followed by
#    more body content cleanly swaps its contents in place (doing the
#    same thing at the very end of the body, just before sectPr,
#    leaves behind a stray empty paragraph), and preserves the exact
#    run split used upstream plus the list numbering/style of the
#    anchor paragraph.
# ------------------------------------------------------------------
$tenure = Find-ParagraphByText $d "Churn and Tenure"

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'
$pPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr>'

$dataUsageBody = $pPr + '<w:r><w:t>Churn and Data Usage</w:t></w:r><w:r><w:t>: The association between the length of time a customer has been with the service (tenure) and their likelihood to churn is analysed</w:t></w:r><w:r><w:t>.</w:t></w:r>'
$statesBody = $pPr + '<w:r><w:t xml:space="preserve">Churn and </w:t></w:r><w:r><w:t>States:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Impact of location on churning ' + [char]0x2013 + ' focus on different states.</w:t></w:r>'
$citiesBody = $pPr + '<w:r><w:t xml:space="preserve">Churn and </w:t></w:r><w:r><w:t>Cities:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Impact of location on churning - focus on different cities.</w:t></w:r>'

$target = $tenure.Range
$anchor = $d.Range($target.End, $target.End)
$anchor.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($tenure.Index + 1)
$newPara.Range.InsertXML($pkgHeader + '<w:body><w:p>' + $dataUsageBody + '</w:p></w:body>' + $pkgFooter)

$dataUsagePara = $d.Paragraphs.Item($tenure.Index + 1)
$anchor2 = $d.Range($dataUsagePara.Range.End, $dataUsagePara.Range.End)
$anchor2.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($dataUsagePara.Index + 1)
$newPara2.Range.InsertXML($pkgHeader + '<w:body><w:p>' + $statesBody + '</w:p></w:body>' + $pkgFooter)

$statesPara = $d.Paragraphs.Item($dataUsagePara.Index + 1)
$anchor3 = $d.Range($statesPara.Range.End, $statesPara.Range.End)
$anchor3.InsertParagraphAfter()
$newPara3 = $d.Paragraphs.Item($statesPara.Index + 1)
$newPara3.Range.InsertXML($pkgHeader + '<w:body><w:p>' + $citiesBody + '</w:p></w:body>' + $pkgFooter)

# ------------------------------------------------------------------
# 2. Now remove the "Advanced Analysis" Heading1 paragraph and the
#    two bulleted paragraphs that followed it (numId=9 list items) -
#    their textual content has been superseded by the three new
#    numId=8 bullets inserted above.
# ------------------------------------------------------------------
$advanced = Find-ParagraphByText $d "Advanced Analysis"
$byRegion = Find-ParagraphByText $d "Churn by Region"
$svcUsage = Find-ParagraphByText $d "Churn and Service Usage"

$svcUsage.Range.Delete()
$byRegion.Range.Delete()
$advanced.Range.Delete()

# ------------------------------------------------------------------
# 3. Mark the built-in "Default Paragraph Font" character style as
#    semi-hidden in the style gallery (matches the styles.xml diff).
# ------------------------------------------------------------------
$dpf = $d.Styles.Item("Default Paragraph Font")
$dpf.SemiHidden = $true
